# Reduce total expenses to better represent a logical income-to-expense
# ratio: the "Rent" expense entries drop from 1200 -> 800 (4 occurrences
# in the transaction log), and the data table gets a white fill applied
# (new solid white fill added to the style table, mirroring the source
# diff's new <fill> entry).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value edits: Rent expense amounts 1200 -> 800 -------------------
$ws.Range("C3").Value  = 800
$ws.Range("C13").Value = 800
$ws.Range("C23").Value = 800
$ws.Range("C33").Value = 800

# --- Formatting: apply a solid white fill across the populated table ------
# (adds a new fill to the workbook's style table, as in the source edit)
$ws.UsedRange.Interior.Color = 16777215
